$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pregunta1"
$ws.Range("B2").Value = "(100, 50)"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "Pregunta1"
$ws.Range("B3").Value = "(200, 80)"
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = "Pregunta2"
$ws.Range("B4").Value = "(150, 60)"
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "Pregunta2"
$ws.Range("B5").Value = "(250, 90)"
$ws.Range("C5").Value = 1
